$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("regional v02")

$ws.Range("B258").Value = "RUFXBRR"
$ws.Range("B259").Value = "RUFXBRF, RUFXBRZ, RUFXBRH"
$ws.Range("C259").Value = "special. Create last step to add special edges"
$ws.Range("B260").Value = "RUAULOGR"
$ws.Range("B261").Value = "RUFXTRF"
$ws.Range("B262").Value = "RUINTRF, RUINTRH"
$ws.Range("C262").Value = "special. Create last step to add special edges"
$ws.Range("B263").Value = "RUPARMF, RUPARMH"
$ws.Range("C263").Value = "special. Create last step to add special edges"
$ws.Range("B264").Value = "RUAULOGR"
$ws.Range("B265").Value = "RUPARRF"
$ws.Range("B266").Value = "RUTRLOGR"
$ws.Range("B267").Value = "RUR002P, RUR003P, RUR009P, RUR012P, RUR013P"
$ws.Range("C267").Value = "special. Create last step to add special edges"
$ws.Range("B268").Value = "RUFXBRR"
$ws.Range("B269").Value = "RUSROBR"
$ws.Range("B270").Value = "RUSROSP"
$ws.Range("B271").Value = "RUTRLOGR"
$ws.Range("B272").Value = "RUAV"
$ws.Range("B273").Value = "RUFU"
$ws.Range("B274").Value = "RUSG"
$ws.Range("B275").Value = "RUSROBUF"
$ws.Range("B276").Value = "RUTRLOGR"
$ws.Range("B277").Value = "RUBU"
$ws.Range("B278").Value = "RUSROHIE"
$ws.Range("B279").Value = "RUTRLOGR"
$ws.Range("B280").Value = "RUSROBR"
$ws.Range("B281").Value = "RUSROSP"
$ws.Range("B282").Value = "RUTRLOGR"
$ws.Range("B283").Value = "RUSROSW"
$ws.Range("B284").Value = "RUTRLOGR"
$ws.Range("B285").Value = "RUSPFXR"
$ws.Range("B286").Value = "RUFXBRR"
$ws.Range("B287").Value = "RUFXTRR"
$ws.Range("B288").Value = "RUPARMR"
$ws.Range("B289").Value = "RUFXBRR"
$ws.Range("B290").Value = "RUSPFXF, RUSPFXZ"
$ws.Range("C290").Value = "special. Create last step to add special edges"
$ws.Range("B291").Value = "RUAULOGR"
$ws.Range("B292").Value = "RUSPGLF, RUSPGLH"
$ws.Range("C292").Value = "special. Create last step to add special edges"
$ws.Range("B293").Value = "RUSPOTF, RUSPOTH"
$ws.Range("C293").Value = "special. Create last step to add special edges"
$ws.Range("B294").Value = "RUFXBRR"
$ws.Range("B295").Value = "RUPARMR"
$ws.Range("B296").Value = "RUBUFFR"
$ws.Range("B297").Value = "RUPARMR"
$ws.Range("B298").Value = "RUBUFFR"
$ws.Range("B299").Value = "RUCCYHR"
$ws.Range("B300").Value = "RUINTRR"
$ws.Range("B301").Value = "RUPARMR"
$ws.Range("B302").Value = "PSRDTEC3"
$ws.Range("B303").Value = "RU702R"
$ws.Range("B304").Value = "RUFXBRR"
$ws.Range("B305").Value = "RUPARMR"
$ws.Range("B306").Value = "RUSPFXR"
$ws.Range("B307").Value = "RUSROHIE"
$ws.Range("B308").Value = "TP8602R"
$ws.Range("B309").Value = "RUSWAPF, RUSWAPH"
$ws.Range("C309").Value = "special. Create last step to add special edges"
$ws.Range("B310").Value = "RUSWGLF, RUSWGLH"
$ws.Range("C310").Value = "special. Create last step to add special edges"
$ws.Range("B311").Value = "RUTENRF, RUTENRH"
$ws.Range("C311").Value = "special. Create last step to add special edges"
$ws.Range("B312").Value = "RUTRLOG"
$ws.Range("B313").Value = "CEELOCT"
$ws.Range("B314").Value = "CEELOCT"

$ws.Range("C312").Select()
